$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (rotated content) ---
$ws.Range("A4").Value = 130937843
$ws.Range("B4").Value = 57884
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("M4").Value = "färska spår"
$ws.Range("Q4").Value = 489760
$ws.Range("R4").Value = 7004232
$ws.Range("AC4").Value = "Ringhack, färska och äldre, i riklig mängd längs flera meter högt upp på en granstam med spår av rikligt sav/kådaflöde."
$ws.Range("AH4").Value = "Granskog"
$ws.Range("AJ4").Value = "gran"
$ws.Range("AK4").Value = "Picea abies"
$ws.Range("AM4").Value = "Trädstam på levande träd"
$ws.Range("AO4").Value = "Stem on living tree # Picea abies"

# --- Row 5 (rotated content) ---
$ws.Range("A5").Value = 130937852
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 489520
$ws.Range("R5").Value = 7004161
$ws.Range("AC5").Value = "Ringhack, äldre, ytliga enstaka längs flera meter på en granstam vid kanten mot yngre skog."
$ws.Range("AH5").Value = "Granskog"
$ws.Range("AJ5").Value = "gran"
$ws.Range("AK5").Value = "Picea abies"
$ws.Range("AM5").Value = "Trädstam på levande träd"
$ws.Range("AO5").Value = "Stem on living tree # Picea abies"

# --- Row 6 (rotated content) ---
$ws.Range("A6").Value = 130937857
$ws.Range("B6").Value = 97878
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221945
$ws.Range("F6").Value = "Revlummer"
$ws.Range("G6").Value = "Lycopodium annotinum"
$ws.Range("H6").Value = "L."
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 489680
$ws.Range("R6").Value = 7004154
$ws.Range("AC6").ClearContents()
$ws.Range("AH6").Value = "Granskog"
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AM6").ClearContents()
$ws.Range("AO6").ClearContents()

# --- Row 7 (rotated content) ---
$ws.Range("A7").Value = 130937863
$ws.Range("B7").Value = 99013
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("I7").Value = "'8"
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("K7").Value = "fullt utvecklade blad"
$ws.Range("Q7").Value = 489799
$ws.Range("R7").Value = 7004245
$ws.Range("AC7").Value = "Minst 8 plantor inom ca 1 m2 yta. Grävdes varsamt fram under snötäcket. Det finns sannolikt betydligt mer knärot på fyndplatsen och i skogsbeståndet där fyndplatsen ligger."
$ws.Range("AH7").Value = "Barrskog"
